$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Convert B2, B9 from numeric 398729 to text shared-string "398729"
$ws.Range("B2").Value = "398729"
$ws.Range("B9").Value = "398729"

# Convert B6, B7 from numeric 372152/396109 to text shared-string "372152"
$ws.Range("B6").Value = "372152"
$ws.Range("B7").Value = "372152"

# Shift D7/E7 dates forward by 9 days
$ws.Range("D7").Value = 43872.206388888888
$ws.Range("E7").Value = 43872.241539351853

# Shift D9/E9 dates back by 10 days
$ws.Range("D9").Value = 43845.011956018519
$ws.Range("E9").Value = 43847.222569444442

# Move the active selection from E7 to E8
$ws.Range("E8").Select()
